$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: C57 and C372 input corrections (cascade through B via cumulative formula) ---
$ws.Range("C57").Value = 8
$ws.Range("C372").Value = 59

# --- Section 2: E/F/G corrections for rows 490-525 ---
$ws.Range("E490").Value = 4
$ws.Range("G491").Value = 3
$ws.Range("E492").Value = 4
$ws.Range("G492").Value = 3
$ws.Range("E493").Value = 2
$ws.Range("G493").Value = 5
$ws.Range("G494").Value = 5
$ws.Range("G495").Value = 5
$ws.Range("G496").Value = 6
$ws.Range("G497").Value = 4
$ws.Range("G499").Value = 4
$ws.Range("G500").Value = 4
$ws.Range("E501").Value = 3
$ws.Range("G501").Value = 4
$ws.Range("E503").Value = 2
$ws.Range("F503").Value = 1
$ws.Range("G503").Value = 4
$ws.Range("G504").Value = 3
$ws.Range("G505").Value = 3
$ws.Range("E506").Value = 1
$ws.Range("F506").Value = 0
$ws.Range("G506").Value = 3
$ws.Range("E507").Value = 2
$ws.Range("F507").Value = 1
$ws.Range("G507").Value = 2
$ws.Range("G508").Value = 2
$ws.Range("G509").Value = 2
$ws.Range("G510").Value = 2
$ws.Range("F511").Value = 1
$ws.Range("G511").Value = 2
$ws.Range("E512").Value = 2
$ws.Range("F512").Value = 2
$ws.Range("E513").Value = 2
$ws.Range("E514").Value = 2
$ws.Range("E517").Value = 2
$ws.Range("F517").Value = 2
$ws.Range("G517").Value = 2
$ws.Range("G518").Value = 2
$ws.Range("E519").Value = 2
$ws.Range("G519").Value = 3
$ws.Range("G520").Value = 3
$ws.Range("G521").Value = 3
$ws.Range("E522").Value = 2
$ws.Range("F522").Value = 2
$ws.Range("E523").Value = 2
$ws.Range("F523").Value = 2
$ws.Range("E524").Value = 2
$ws.Range("E525").Value = 2

# --- Section 3: C/E/F/G(/L/M) corrections & new data for rows 526-531 ---
$ws.Range("C526").Value = 26
$ws.Range("E526").Value = 2
$ws.Range("G526").Value = 3
$ws.Range("C527").Value = 29
$ws.Range("E527").Value = 2
$ws.Range("F527").Value = 2
$ws.Range("G527").Value = 4
$ws.Range("C528").Value = 27
$ws.Range("E528").Value = 2
$ws.Range("F528").Value = 2
$ws.Range("G528").Value = 3
$ws.Range("L528").Value = 0
$ws.Range("M528").Value = 0
$ws.Range("C529").Value = 20
$ws.Range("E529").Value = 2
$ws.Range("F529").Value = 2
$ws.Range("G529").Value = 4
$ws.Range("L529").Value = 0
$ws.Range("M529").Value = 0
$ws.Range("C530").Value = 25
$ws.Range("E530").Value = 2
$ws.Range("F530").Value = 2
$ws.Range("G530").Value = 5
$ws.Range("L530").Value = 0
$ws.Range("M530").Value = 0
$ws.Range("C531").Value = 2
$ws.Range("E531").Value = 2
$ws.Range("F531").Value = 2
$ws.Range("G531").Value = 5
$ws.Range("L531").Value = 0
$ws.Range("M531").Value = 0

# --- Section 4: restore the active cell selection as recorded in the saved view state ---
$ws.Range("O524").Select()
